# Auto-generated edit script: updates Leve profit-calculation cells (columns H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 106.166664
$ws.Range("I38").Value = 106.166664
$ws.Range("K38").Value = 318.499992
$ws.Range("M38").Value = 53.50000799999998
$ws.Range("H106").Value = 1984.7273
$ws.Range("I106").Value = 2889.3333
$ws.Range("K106").Value = 2889.3333
$ws.Range("M106").Value = -2258.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 12503031
$ws.Range("I45").Value = 3187.4285
$ws.Range("K45").Value = 3187.4285
$ws.Range("M45").Value = -2810.4285
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("I110").Value = 1466.6666
$ws.Range("J110").Value = 1200
$ws.Range("K110").Value = 1466.6666
$ws.Range("L110").Value = 1200
$ws.Range("M110").Value = 578.3334
$ws.Range("N110").Value = -5290
$ws.Range("H122").Value = 2013.625
$ws.Range("I122").Value = 1881.7
$ws.Range("K122").Value = 5645.1
$ws.Range("M122").Value = -3195.1
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()
$ws.Range("H139").Value = 149714.5
$ws.Range("J139").Value = 149714.5
$ws.Range("L139").Value = 149714.5
$ws.Range("N139").Value = -159994.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 586616.5
$ws.Range("I20").Value = 586616.5
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 586616.5
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -586369.5
$ws.Range("N20").ClearContents()
$ws.Range("H64").Value = 2449.4
$ws.Range("I64").Value = 1750
$ws.Range("J64").Value = 2624.25
$ws.Range("K64").Value = 1750
$ws.Range("L64").Value = 2624.25
$ws.Range("M64").Value = -1525
$ws.Range("N64").Value = -3074.25
$ws.Range("H67").Value = 2449.4
$ws.Range("I67").Value = 1750
$ws.Range("J67").Value = 2624.25
$ws.Range("K67").Value = 1750
$ws.Range("L67").Value = 2624.25
$ws.Range("M67").Value = -970
$ws.Range("N67").Value = -4184.25
$ws.Range("H80").Value = 1130.125
$ws.Range("I80").Value = 867.3333
$ws.Range("J80").Value = 1287.8
$ws.Range("K80").Value = 867.3333
$ws.Range("L80").Value = 1287.8
$ws.Range("M80").Value = 130.6667
$ws.Range("N80").Value = -3283.8
$ws.Range("H83").Value = 1130.125
$ws.Range("I83").Value = 867.3333
$ws.Range("J83").Value = 1287.8
$ws.Range("K83").Value = 4336.6665
$ws.Range("L83").Value = 6439
$ws.Range("M83").Value = 655.3334999999997
$ws.Range("N83").Value = -16423
$ws.Range("H87").Value = 75000
$ws.Range("J87").Value = 75000
$ws.Range("L87").Value = 75000
$ws.Range("N87").Value = -77496
$ws.Range("H90").Value = 75000
$ws.Range("J90").Value = 75000
$ws.Range("L90").Value = 225000
$ws.Range("N90").Value = -237480
$ws.Range("H134").Value = 3292.6829
$ws.Range("I134").Value = 2662.162
$ws.Range("K134").Value = 7986.485999999999
$ws.Range("M134").Value = -5451.485999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 775
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 974.6667
$ws.Range("I22").Value = 1049.6
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 1049.6
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -699.5999999999999
$ws.Range("N22").Value = -1300
$ws.Range("H31").Value = 2556.04
$ws.Range("J31").Value = 4379.4
$ws.Range("L31").Value = 4379.4
$ws.Range("N31").Value = -4969.4
$ws.Range("H34").Value = 2556.04
$ws.Range("J34").Value = 4379.4
$ws.Range("L34").Value = 4379.4
$ws.Range("N34").Value = -4783.4
$ws.Range("H86").Value = 6708.25
$ws.Range("J86").Value = 7757
$ws.Range("L86").Value = 7757
$ws.Range("N86").Value = -10003
$ws.Range("H89").Value = 6708.25
$ws.Range("J89").Value = 7757
$ws.Range("L89").Value = 38785
$ws.Range("N89").Value = -50017
$ws.Range("H94").Value = 1858.3334
$ws.Range("I94").Value = 1787.5
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 1787.5
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -1336.5
$ws.Range("N94").Value = -2902
$ws.Range("H107").Value = 887
$ws.Range("I107").Value = 1048.75
$ws.Range("K107").Value = 1048.75
$ws.Range("M107").Value = 871.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 1000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2468
$ws.Range("N45").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H88").Value = 4674
$ws.Range("J88").Value = 4674
$ws.Range("L88").Value = 14022
$ws.Range("N88").Value = -14878
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H91").Value = 4674
$ws.Range("J91").Value = 4674
$ws.Range("L91").Value = 14022
$ws.Range("N91").Value = -16986

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 213.77777
$ws.Range("I2").Value = 7.5
$ws.Range("K2").Value = 7.5
$ws.Range("M2").Value = 105.5
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H59").Value = 9411.111000000001
$ws.Range("I59").Value = 9411.111000000001
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 9411.111000000001
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -8828.111000000001
$ws.Range("N59").ClearContents()
$ws.Range("H70").Value = 171000
$ws.Range("I70").Value = 171000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 171000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -170730
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 171000
$ws.Range("I73").Value = 171000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 171000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -170064
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 1756.25
$ws.Range("I80").Value = 1849
$ws.Range("J80").Value = 1725.3334
$ws.Range("K80").Value = 1849
$ws.Range("L80").Value = 1725.3334
$ws.Range("M80").Value = -851
$ws.Range("N80").Value = -3721.3334
$ws.Range("H83").Value = 1756.25
$ws.Range("I83").Value = 1849
$ws.Range("J83").Value = 1725.3334
$ws.Range("K83").Value = 9245
$ws.Range("L83").Value = 8626.666999999999
$ws.Range("M83").Value = -4253
$ws.Range("N83").Value = -18610.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4178
$ws.Range("I7").Value = 2999.5
$ws.Range("J7").Value = 5524.857
$ws.Range("K7").Value = 2999.5
$ws.Range("L7").Value = 5524.857
$ws.Range("M7").Value = -2887.5
$ws.Range("N7").Value = -5748.857
$ws.Range("H22").Value = 1085
$ws.Range("I22").Value = 1396
$ws.Range("J22").Value = 968.375
$ws.Range("K22").Value = 1396
$ws.Range("L22").Value = 968.375
$ws.Range("M22").Value = -1101
$ws.Range("N22").Value = -1558.375
$ws.Range("H27").Value = 1085
$ws.Range("I27").Value = 1396
$ws.Range("J27").Value = 968.375
$ws.Range("K27").Value = 1396
$ws.Range("L27").Value = 968.375
$ws.Range("M27").Value = -1289
$ws.Range("N27").Value = -1182.375
$ws.Range("H33").Value = 14799.5
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H46").Value = 2234.7144
$ws.Range("I46").Value = 1425
$ws.Range("J46").Value = 2558.6
$ws.Range("K46").Value = 1425
$ws.Range("L46").Value = 2558.6
$ws.Range("M46").Value = -1237
$ws.Range("N46").Value = -2934.6
$ws.Range("H126").Value = 4178
$ws.Range("I126").Value = 2999.5
$ws.Range("J126").Value = 5524.857
$ws.Range("K126").Value = 8998.5
$ws.Range("L126").Value = 16574.571
$ws.Range("M126").Value = -6528.5
$ws.Range("N126").Value = -21514.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 17798.572
$ws.Range("I107").Value = 37297
$ws.Range("K107").Value = 111891
$ws.Range("M107").Value = -109971
$ws.Range("H132").Value = 1502.4
$ws.Range("I132").Value = 1261.0834
$ws.Range("J132").Value = 2467.6667
$ws.Range("K132").Value = 3783.2502
$ws.Range("L132").Value = 7403.000100000001
$ws.Range("M132").Value = -1253.2502
$ws.Range("N132").Value = -12463.0001

Write-Host "Updated Leve profit figures across all job sheets."
